$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ALC_refs = @("H9", "I9", "K9", "M9", "H11", "I11", "K11", "M11", "H18", "I18", "J18", "K18", "L18", "M18", "H19", "J19", "L19", "N19", "H32", "I32", "J32", "K32", "L32", "M32", "N32", "H38", "I38", "K38", "M38", "H64", "J64", "L64", "N64", "H67", "J67", "L67", "N67", "H74", "I74", "K74", "M74", "H77", "I77", "K77", "M77", "H88", "I88", "J88", "K88", "L88", "N88", "H91", "I91", "J91", "K91", "L91", "N91", "H100", "I100", "K100", "M100", "H107", "I107", "J107", "K107", "L107", "M107", "N107", "H112", "J112", "L112", "N112", "H125", "I125", "J125", "K125", "L125", "M125", "N125", "H131", "I131", "J131", "K131", "L131", "M131", "N131", "H132", "I132", "J132", "K132", "L132", "M132", "N132", "H137", "I137", "J137", "K137", "L137", "M137", "N137", "H138", "J138", "L138", "N138")
$ALC_vals = @(50000716, 125000140, 125000140, -124999971, 67.30768999999999, 67.30768999999999, 67.30768999999999, 72.69231000000001, 5159.467, 5159.467, 0, 5159.467, 0, -4875.467, 850.9259, 1039.8334, 1039.8334, -1389.8334, 3741.158, 3676.3333, 3799.5, 3676.3333, 3799.5, -3350.3333, -4451.5, 2257.7273, 1041.875, 3125.625, -2753.625, 27182.27, 5720.1, 5720.1, -6216.1, 27182.27, 5720.1, 5720.1, -7436.1, 5055.857, 3398.2, 3398.2, -2462.2, 5055.857, 3398.2, 16991, -12311, 6383.25, 0, 6383.25, 0, 6383.25, -7195.25, 6383.25, 0, 6383.25, 0, 6383.25, -9191.25, 9489180, 10876848, 10876848, -10876307, 9662.416999999999, 13677.375, 1632.5, 13677.375, 1632.5, -11757.375, -5472.5, 31546.523, 31546.523, 94639.569, -96855.569, 1499.5834, 334, 4996.3335, 3006, 44967.0015, -546, -49887.0015, 4166.722, 2085.7144, 5491, 6257.1432, 16473, -1217.1432, -26553, 4809.04, 4853.5654, 4297, 14560.6962, 12891, -12030.6962, -17951, 337081.8, 1146358, 6014.273, 3439074, 18042.819, -3436524, -23142.819, 204819.48, 6839.75, 20519.25, -30799.25)
for ($i = 0; $i -lt $ALC_refs.Length; $i++) {
    $ws.Range($ALC_refs[$i]).Value = $ALC_vals[$i]
}
$ALC_clear = @("N18", "M88", "M91")
foreach ($ref in $ALC_clear) {
    $ws.Range($ref).ClearContents()
}

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ARM_refs = @("H2", "I2", "K2", "M2", "H10", "J10", "L10", "N10", "H32", "I32", "K32", "M32", "H45", "I45", "J45", "K45", "L45", "M45", "N45", "H61", "I61", "J61", "K61", "L61", "M61", "N61", "H63", "I63", "J63", "K63", "L63", "M63", "H66", "I66", "J66", "K66", "L66", "M66", "H74", "I74", "K74", "M74", "H77", "I77", "K77", "M77", "H110", "I110", "J110", "K110", "L110", "M110", "N110", "H116", "I116", "K116", "M116", "H132", "I132", "J132", "K132", "L132", "M132", "N132", "H136", "I136", "J136", "K136", "L136", "M136", "N136", "H139", "J139", "L139", "N139")
$ARM_vals = @(94541.664, 17500, 17500, -17387, 22998.5, 3, 3, -343, 2686.717, 2752.8628, 2752.8628, -2465.8628, 156943.42, 270953.62, 4929.8335, 270953.62, 4929.8335, -270576.62, -5683.8335, 5753.2856, 6384.8184, 3437.6667, 6384.8184, 3437.6667, -6172.8184, -3861.6667, 8345.888999999999, 8345.888999999999, 0, 8345.888999999999, 0, -7659.888999999999, 8345.888999999999, 8345.888999999999, 0, 41729.44499999999, 0, -38297.44499999999, 2846.2666, 2482.8333, 2482.8333, -1608.8333, 2846.2666, 2482.8333, 12414.1665, -8046.166499999999, 1839.8572, 1813.1666, 2000, 1813.1666, 2000, 231.8334, -6090, 94541.664, 17500, 17500, -15206, 3011.7917, 2312.946, 5362.4546, 6938.838, 16087.3638, -4408.838, -21147.3638, 5753.2856, 6384.8184, 3437.6667, 19154.4552, 10313.0001, -16604.4552, -15413.0001, 116444.4, 116444.4, 116444.4, -126724.4)
for ($i = 0; $i -lt $ARM_refs.Length; $i++) {
    $ws.Range($ARM_refs[$i]).Value = $ARM_vals[$i]
}
$ARM_clear = @("N63", "N66")
foreach ($ref in $ARM_clear) {
    $ws.Range($ref).ClearContents()
}

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$BSM_refs = @("H3", "I3", "K3", "M3", "H20", "I20", "J20", "K20", "L20", "M20", "N20", "H86", "I86", "J86", "K86", "L86", "M86", "N86", "H89", "I89", "J89", "K89", "L89", "M89", "N89", "H105", "I105", "J105", "K105", "L105", "M105", "N105", "H107", "I107", "J107", "K107", "L107", "M107", "N107", "H134", "I134", "K134", "M134")
$BSM_vals = @(94541.664, 17500, 17500, -17386, 5218.8887, 1874, 9400, 1874, 9400, -1627, -9894, 4764.1035, 5136.654, 1535.3334, 5136.654, 1535.3334, -4013.654, -3781.3334, 4764.1035, 5136.654, 1535.3334, 25683.27, 7676.666999999999, -20067.27, -18908.667, 82606.234, 116653.664, 5999.5, 116653.664, 5999.5, -114906.664, -9493.5, 3737.16, 3464.1765, 4317.25, 3464.1765, 4317.25, -1544.1765, -8157.25, 5737.75, 5921.946, 17765.838, -15230.838)
for ($i = 0; $i -lt $BSM_refs.Length; $i++) {
    $ws.Range($BSM_refs[$i]).Value = $BSM_vals[$i]
}

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$CRP_refs = @("H16", "I16", "J16", "K16", "L16", "M16", "N16", "H31", "I31", "J31", "K31", "L31", "M31", "N31", "H34", "I34", "J34", "K34", "L34", "M34", "N34", "H36", "I36", "K36", "M36", "H40", "I40", "K40", "M40", "H58", "I58", "J58", "K58", "L58", "M58", "N58", "H62", "I62", "J62", "K62", "L62", "M62", "N62", "H65", "I65", "J65", "K65", "L65", "M65", "N65", "H94", "J94", "L94", "N94", "H105", "I105", "K105", "M105", "H113", "I113", "J113", "K113", "L113", "M113", "N113", "H132", "I132", "K132", "M132", "H134", "I134", "K134", "M134", "H136", "I136", "J136", "K136", "L136", "M136", "N136", "H140", "J140", "L140", "N140")
$CRP_vals = @(2365.4546, 2113.3333, 3500, 2113.3333, 3500, -1826.3333, -4074, 2547.776, 1832.1666, 2630.3462, 1832.1666, 2630.3462, -1537.1666, -3220.3462, 2547.776, 1832.1666, 2630.3462, 1832.1666, 2630.3462, -1630.1666, -3034.3462, 16682.666, 9998, 9998, -9610, 16682.666, 9998, 9998, -9838, 7150.5557, 7675.087, 6222.5386, 7675.087, 6222.5386, -7472.087, -6628.5386, 9626.857, 8272.286, 10304.143, 8272.286, 10304.143, -7648.286, -11552.143, 9626.857, 8272.286, 10304.143, 41361.43, 51520.715, -38241.43, -57760.715, 2237.9412, 1385.25, 1385.25, -2287.25, 163685.16, 235294.22, 235294.22, -233547.22, 2365.4546, 2113.3333, 3500, 2113.3333, 3500, 56.66670000000022, -7840, 17685.055, 7921.7, 23765.1, -21235.1, 4176050.2, 5219610.5, 15658831.5, -15656296.5, 7150.5557, 7675.087, 6222.5386, 23025.261, 18667.6158, -20475.261, -23767.6158, 79104.336, 79104.336, 79104.336, -89464.336)
for ($i = 0; $i -lt $CRP_refs.Length; $i++) {
    $ws.Range($CRP_refs[$i]).Value = $CRP_vals[$i]
}

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$CUL_refs = @("H5", "I5", "K5", "M5", "H7", "I7", "K7", "M7", "H56", "I56", "K56", "M56", "H62", "I62", "K62", "M62", "H65", "I65", "K65", "M65", "H92", "J92", "L92", "N92", "H107", "I107", "J107", "K107", "L107", "M107", "N107", "H115", "I115", "J115", "K115", "L115", "M115", "N115", "H121", "J121", "L121", "N121", "H132", "I132", "J132", "K132", "L132", "M132", "N132", "H135", "I135", "K135", "M135", "H137", "J137", "L137", "N137")
$CUL_vals = @(527726.9, 911.7143, 2735.1429, -2623.1429, 2500000, 2500000, 7500000, -7499888, 5441.5, 5441.5, 5441.5, -4911.5, 300, 300, 900, -214, 300, 300, 2700, 732, 1465.6666, 1465.6666, 4396.9998, -6892.9998, 2537.5957, 880.3570999999999, 3240.6667, 2641.0713, 9722.000100000001, -721.0712999999996, -13562.0001, 2613.9, 312.4, 4915.4, 937.1999999999999, 14746.2, 237.8000000000001, -17096.2, 2177951, 2223130.5, 6669391.5, -6672011.5, 11147406, 1750, 12862122, 15750, 115759098, -13220, -115764158, 527726.9, 911.7143, 8205.4287, -5670.4287, 8813.77, 9008.272000000001, 27024.816, -37224.81600000001)
for ($i = 0; $i -lt $CUL_refs.Length; $i++) {
    $ws.Range($CUL_refs[$i]).Value = $CUL_vals[$i]
}

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$GSM_refs = @("H70", "I70", "J70", "K70", "L70", "M70", "N70", "H73", "I73", "J73", "K73", "L73", "M73", "N73", "H80", "I80", "J80", "K80", "L80", "M80", "N80", "H83", "I83", "J83", "K83", "L83", "M83", "N83", "H102", "I102", "K102", "M102", "H104", "J104", "L104", "H122", "I122", "J122", "K122", "L122", "M122", "N122", "H126", "J126", "L126", "N126", "H132", "I132", "J132", "K132", "L132", "M132", "N132")
$GSM_vals = @(8613.916999999999, 5858.875, 14124, 5858.875, 14124, -5588.875, -14664, 8613.916999999999, 5858.875, 14124, 5858.875, 14124, -4922.875, -15996, 6168.75, 8531.143, 2861.4, 8531.143, 2861.4, -7533.143, -4857.4, 6168.75, 8531.143, 2861.4, 42655.715, 14307, -37663.715, -24291, 11058.625, 12340.571, 12340.571, -10718.571, 0, 0, 0, 13661.286, 16196.363, 4366, 48589.089, 13098, -46139.089, -17998, 17904.945, 12365.833, 37097.499, -42037.499, 2657, 2102.303, 6318, 6306.909, 18954, -3776.909, -24014)
for ($i = 0; $i -lt $GSM_refs.Length; $i++) {
    $ws.Range($GSM_refs[$i]).Value = $GSM_vals[$i]
}
$GSM_clear = @("N104")
foreach ($ref in $GSM_clear) {
    $ws.Range($ref).ClearContents()
}

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$LTW_refs = @("H7", "J7", "L7", "N7", "H16", "I16", "K16", "M16", "H38", "I38", "J38", "K38", "L38", "M38", "N38", "H40", "I40", "J40", "K40", "L40", "M40", "N40", "H42", "I42", "K42", "M42", "H49", "I49", "K49", "M49", "H61", "I61", "J61", "K61", "L61", "M61", "N61", "H100", "I100", "K100", "M100", "H113", "I113", "J113", "K113", "L113", "M113", "N113", "H122", "I122", "K122", "M122", "H126", "J126", "L126", "N126", "H132", "I132", "J132", "K132", "L132", "M132", "N132", "H136", "I136", "J136", "K136", "L136", "M136", "N136", "H141", "J141", "L141", "N141")
$LTW_vals = @(64390.57, 8750, 8750, -8974, 1973.2, 1973.2, 1973.2, -1803.2, 39717.43, 14030, 43998.668, 14030, 43998.668, -13620, -44818.668, 87108.78, 119664.836, 21996.666, 119664.836, 21996.666, -119528.836, -22268.666, 53341.332, 60012.5, 60012.5, -59449.5, 53341.332, 60012.5, 60012.5, -59865.5, 3680.8, 3544, 4000, 3544, 4000, -3342, -4404, 11983.5, 20467.666, 20467.666, -19926.666, 3680.8, 3544, 4000, 3544, 4000, -1374, -8340, 4672.4443, 4672.4443, 14017.3329, -11567.3329, 64390.57, 8750, 26250, -31190, 2940949.8, 3651187.8, 99998, 10953563.4, 299994, -10951033.4, -305054, 7710.8, 5884.6665, 10450, 17653.9995, 31350, -15103.9995, -36450, 40000, 40000, 40000, -50360)
for ($i = 0; $i -lt $LTW_refs.Length; $i++) {
    $ws.Range($LTW_refs[$i]).Value = $LTW_vals[$i]
}

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$WVR_refs = @("H23", "I23", "K23", "M23", "H101", "J101", "L101", "N101", "H107", "I107", "K107", "M107", "H111", "J111", "L111", "N111", "H122", "J122", "L122", "N122", "H126", "J126", "L126", "N126")
$WVR_vals = @(1478.4, 1478.4, 1478.4, -1249.4, 26250, 26250, 26250, -32740, 33765.9, 4320.4, 12961.2, -11041.2, 60000, 60000, 60000, -68180, 8239.5, 10832.5, 32497.5, -37397.5, 30223.688, 6298, 18894, -23834)
for ($i = 0; $i -lt $WVR_refs.Length; $i++) {
    $ws.Range($WVR_refs[$i]).Value = $WVR_vals[$i]
}
